# Refresh the "cryptos" price/volume snapshot (GitHub Actions style data pull).
# Price (col D) and Volume(1h) (col E) values are plain text in this sheet
# (e.g. "67.564.40", "0.140", "  -0.43%  "), so we write them with a leading
# apostrophe to force text and then strip the resulting quote-prefix style so
# the cell keeps its original (unstyled) look -- this avoids Excel silently
# re-interpreting numeric-looking strings (like "0.140") as numbers and
# dropping significant trailing/leading zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $range = $ws.Range($cell)
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextCell "D2" '67.564.40'
Set-TextCell "E2" '  -0.43%  '
Set-TextCell "D3" '3.524.11'
Set-TextCell "E3" '  -0.66%  '
Set-TextCell "E4" '  -0.11%  '
Set-TextCell "D5" '609.37'
Set-TextCell "E5" '  -1.46%  '
Set-TextCell "D6" '151.39'
Set-TextCell "E6" '  -1.27%  '
Set-TextCell "D7" '3.523.37'
Set-TextCell "E7" '  -0.64%  '
Set-TextCell "E8" '  -0.18%  '
Set-TextCell "E9" '  -0.56%  '
Set-TextCell "D10" '0.140'
Set-TextCell "E10" '  -0.88%  '
Set-TextCell "D11" '7.07'
Set-TextCell "E11" '  +1.76%  '
Set-TextCell "E12" '  -1.43%  '
Set-TextCell "D13" '0.0000220'
Set-TextCell "E13" '  -1.78%  '
Set-TextCell "D14" '4.124.86'
Set-TextCell "E14" '  -0.49%  '
Set-TextCell "D15" '31.97'
Set-TextCell "E15" '  -0.05%  '
Set-TextCell "D16" '3.530.57'
Set-TextCell "E16" '  -0.42%  '
Set-TextCell "D17" '67.491.38'
Set-TextCell "E17" '  -0.58%  '
Set-TextCell "E18" '  +0.19%  '
Set-TextCell "E19" '  +0.35%  '
Set-TextCell "E20" '  -2.33%  '
Set-TextCell "D21" '446.52'
Set-TextCell "E21" '  -2.58%  '
Set-TextCell "D22" '9.32'
Set-TextCell "E22" '  -3.96%  '
Set-TextCell "E23" '  -2.56%  '
Set-TextCell "D24" '77.38'
Set-TextCell "E24" '  -0.38%  '
Set-TextCell "E25" '  +12.13%  '
Set-TextCell "D26" '3.665.44'
Set-TextCell "E26" '  -0.64%  '
Set-TextCell "E27" '  +0.10%  '
Set-TextCell "D28" '10.18'
Set-TextCell "E28" '  -4.26%  '
Set-TextCell "D29" '8.35'
Set-TextCell "E29" '  +0.41%  '
Set-TextCell "E30" '  -2.66%  '
Set-TextCell "E31" '  -3.52%  '
Set-TextCell "D32" '0.999'
Set-TextCell "E32" '  -0.14%  '
Set-TextCell "D33" '0.164'
Set-TextCell "E33" '  +4.63%  '
Set-TextCell "D34" '25.77'
Set-TextCell "E34" '  -0.63%  '
Set-TextCell "D35" '6.14'
Set-TextCell "E35" '  -0.15%  '
Set-TextCell "D36" '3.516.63'
Set-TextCell "E36" '  -0.79%  '
Set-TextCell "E37" '  -3.09%  '
Set-TextCell "E38" '  -0.23%  '
Set-TextCell "E39" '  +0.04%  '
Set-TextCell "D40" '0.999'
Set-TextCell "E40" '  -0.18%  '
Set-TextCell "D41" '177.23'
Set-TextCell "E41" '  -0.29%  '
Set-TextCell "D42" '2.19'
Set-TextCell "E42" '  +3.88%  '
Set-TextCell "E43" '  -0.68%  '
Set-TextCell "D44" '5.45'
Set-TextCell "E44" '  -2.92%  '
Set-TextCell "E45" '  -1.42%  '
Set-TextCell "D46" '45.57'
Set-TextCell "E46" '  -0.66%  '
Set-TextCell "D49" '27.21'
Set-TextCell "E49" '  -4.50%  '
Set-TextCell "E50" '  -1.48%  '
Set-TextCell "D51" '0.997'
Set-TextCell "E51" '  -0.96%  '

# Row 47/48 swap: ONDO <-> dogwifhat, with updated price/volume values
Set-TextCell "B47" 'dogwifhat'
Set-TextCell "C47" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell "D47" '2.62'
Set-TextCell "E47" '  +2.30%  '
Set-TextCell "B48" 'ONDO'
Set-TextCell "C48" 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextCell "D48" '2.62'
Set-TextCell "E48" '  +5.30%  '
